$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- New OpenMP timing values in column H (data rows) ---
$ws.Range("H5").Value = 1784.4763
$ws.Range("H7").Value = 63.197
$ws.Range("H8").Value = 215.2483

# H11 / H12 / H14 are brand-new cells in this block; give them the same
# "0.00" numeric format used throughout column H (matches style s="2").
$ws.Range("H11").NumberFormat = "0.00"
$ws.Range("H11").Font.Bold = $false
$ws.Range("H11").Value = 290.6166

$ws.Range("H12").NumberFormat = "0.00"
$ws.Range("H12").Font.Bold = $false
$ws.Range("H12").Value = 990.417

$ws.Range("H14").NumberFormat = "0.00"
$ws.Range("H14").Font.Bold = $false
$ws.Range("H14").Value = 1771.2599

# --- OpenMP speedup formulas (C / H) added to the summary block ---
$ws.Range("H19").Formula = "=C5/H5"
$ws.Range("H19").NumberFormat = "0.00"
$ws.Range("H19").Font.Bold = $true

$ws.Range("H21").Formula = "=C7/H7"

$ws.Range("H22").Formula = "=C8/H8"
$ws.Range("H22").NumberFormat = "0.00"
$ws.Range("H22").Font.Bold = $true

$ws.Range("H25").Formula = "=C11/H11"
$ws.Range("H25").NumberFormat = "0.00"
$ws.Range("H25").Font.Bold = $true

$ws.Range("H26").Formula = "=C12/H12"
$ws.Range("H26").NumberFormat = "0.00"
$ws.Range("H26").Font.Bold = $true

$ws.Range("H28").Formula = "=C14/H14"
$ws.Range("H28").NumberFormat = "0.00"
$ws.Range("H28").Font.Bold = $false

# --- Formatting touch-ups on existing speedup cells (F/G columns) ---
$ws.Range("F21").NumberFormat = "0.00"
$ws.Range("F21").Font.Bold = $true

$ws.Range("G19").NumberFormat = "0.00"
$ws.Range("G19").Font.Bold = $false

$ws.Range("G22").NumberFormat = "0.00"
$ws.Range("G22").Font.Bold = $false

$ws.Range("G25").NumberFormat = "0.00"
$ws.Range("G25").Font.Bold = $false

$ws.Range("F26").NumberFormat = "0.00"
$ws.Range("F26").Font.Bold = $false

# --- Reposition the selection to match the saved UI state ---
$ws.Range("G30").Select()
